# Use Case Add Property Type.docx - camelCase the field names referenced in
# the "Estate Agent Enters Required Data" / "Validate Data Entered" / "Save
# data..." steps (Property Type Code -> typeCode, Property Type / type
# Description -> typeDescription) and mark them up the way Word's spell
# checker does when it flags a camelCase word (proofErr spellStart/spellEnd
# around a run split on the lower-case prefix). Also stamp a
# lastRenderedPageBreak on the "Alternate Scenarios" heading run.

function Set-ParagraphXml {
    param([string]$OldText, [string]$NewParaXml)

    # Paragraph.Range.Text always carries a trailing paragraph mark (CR,
    # 0x0D); when the paragraph is also the last one in a table cell it
    # carries a further cell-mark (0x07) after that. Compare on the
    # visible text only so both cases match the same lookup string.
    $wanted = $OldText.TrimEnd([char]13, [char]7)

    $doc = $word.ActiveDocument
    $found = $false
    for ($i = 1; $i -le $doc.Paragraphs.Count; $i++) {
        $para = $doc.Paragraphs.Item($i)
        $actual = $para.Range.Text.TrimEnd([char]13, [char]7)
        if ($actual -eq $wanted) {
            $pkg = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" xmlns:xml="http://www.w3.org/XML/1998/namespace"><w:body>' + $NewParaXml + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
            $para.Range.InsertXML($pkg)
            $found = $true
            break
        }
    }
    if (-not $found) {
        Write-Output "NOT FOUND: [$wanted]"
    }
}

# 1) "Property Type Code (2 characters)" bullet under Step 3
$xml1 = @'
<w:p w14:paraId="0B77AFC4" w14:textId="27F743DC" w:rsidR="00325DB6" w:rsidRPr="00325DB6" w:rsidRDefault="00325DB6" w:rsidP="00325DB6"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:rPr><w:bCs/><w:sz w:val="24"/></w:rPr></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:bCs/><w:sz w:val="24"/></w:rPr><w:t>t</w:t></w:r><w:r><w:rPr><w:bCs/><w:sz w:val="24"/></w:rPr><w:t>ypeCode</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidR="00EC54B0"><w:rPr><w:bCs/><w:sz w:val="24"/></w:rPr><w:t xml:space="preserve"> (2 characters)</w:t></w:r></w:p>
'@
Set-ParagraphXml "Property Type Code (2 characters)" $xml1

# 2) "Property Type Description" bullet under Step 3
$xml2 = @'
<w:p w14:paraId="33170D3A" w14:textId="62AE302C" w:rsidR="00325DB6" w:rsidRPr="00325DB6" w:rsidRDefault="00325DB6" w:rsidP="00325DB6"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:rPr><w:bCs/><w:sz w:val="24"/></w:rPr></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:bCs/><w:sz w:val="24"/></w:rPr><w:t>t</w:t></w:r><w:r><w:rPr><w:bCs/><w:sz w:val="24"/></w:rPr><w:t>ypeDescription</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p>
'@
Set-ParagraphXml "Property Type Description" $xml2

# 3) "Type Code must not already exist." bullet under Step 4
$xml3 = @'
<w:p w14:paraId="76490662" w14:textId="22DB25F1" w:rsidR="00325DB6" w:rsidRDefault="00325DB6" w:rsidP="00325DB6"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr><w:rPr><w:bCs/><w:sz w:val="24"/></w:rPr></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:bCs/><w:sz w:val="24"/></w:rPr><w:t>t</w:t></w:r><w:r><w:rPr><w:bCs/><w:sz w:val="24"/></w:rPr><w:t>ypeCode</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:bCs/><w:sz w:val="24"/></w:rPr><w:t xml:space="preserve"> must not already exist.</w:t></w:r></w:p>
'@
Set-ParagraphXml "Type Code must not already exist." $xml3

# 4) "Type Code must be 2 characters in length." bullet under Step 4
$xml4 = @'
<w:p w14:paraId="4C4B4420" w14:textId="65C80ACF" w:rsidR="00EC54B0" w:rsidRPr="00325DB6" w:rsidRDefault="00EC54B0" w:rsidP="00325DB6"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr><w:rPr><w:bCs/><w:sz w:val="24"/></w:rPr></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:bCs/><w:sz w:val="24"/></w:rPr><w:t>t</w:t></w:r><w:r><w:rPr><w:bCs/><w:sz w:val="24"/></w:rPr><w:t>ypeCode</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:bCs/><w:sz w:val="24"/></w:rPr><w:t xml:space="preserve"> must be 2 characters in length.</w:t></w:r></w:p>
'@
Set-ParagraphXml "Type Code must be 2 characters in length." $xml4

# 5) "Property type Description must not be empty." bullet under Step 4
$xml5 = @'
<w:p w14:paraId="3139207B" w14:textId="5E8B58D1" w:rsidR="00EC54B0" w:rsidRPr="00EC54B0" w:rsidRDefault="00325DB6" w:rsidP="00EC54B0"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr><w:rPr><w:bCs/><w:sz w:val="24"/></w:rPr></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:bCs/><w:sz w:val="24"/></w:rPr><w:t>typeDescription</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:bCs/><w:sz w:val="24"/></w:rPr><w:t xml:space="preserve"> must not be empty.</w:t></w:r></w:p>
'@
Set-ParagraphXml "Property type Description must not be empty." $xml5

# 6) "Property Type Code." bullet under Step 5
$xml6 = @'
<w:p w14:paraId="032D13E7" w14:textId="2C647A55" w:rsidR="00A0340D" w:rsidRPr="00A0340D" w:rsidRDefault="00A0340D" w:rsidP="00A0340D"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="3"/></w:numPr><w:rPr><w:bCs/><w:sz w:val="24"/></w:rPr></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:bCs/><w:sz w:val="24"/></w:rPr><w:t>t</w:t></w:r><w:r><w:rPr><w:bCs/><w:sz w:val="24"/></w:rPr><w:t>ypeCode</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:bCs/><w:sz w:val="24"/></w:rPr><w:t>.</w:t></w:r></w:p>
'@
Set-ParagraphXml "Property Type Code." $xml6

# 7) "Property Type Description." bullet under Step 5
$xml7 = @'
<w:p w14:paraId="0A5C60AC" w14:textId="2F311B0E" w:rsidR="00A0340D" w:rsidRDefault="00A0340D" w:rsidP="00A0340D"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="3"/></w:numPr><w:rPr><w:bCs/><w:sz w:val="24"/></w:rPr></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:bCs/><w:sz w:val="24"/></w:rPr><w:t>t</w:t></w:r><w:r><w:rPr><w:bCs/><w:sz w:val="24"/></w:rPr><w:t>ypeDescription</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:bCs/><w:sz w:val="24"/></w:rPr><w:t>.</w:t></w:r></w:p>
'@
Set-ParagraphXml "Property Type Description." $xml7

# 8) "Alternate Scenarios" heading - add lastRenderedPageBreak before the text
$xml8 = @'
<w:p w14:paraId="6810B955" w14:textId="77777777" w:rsidR="00BB2E43" w:rsidRPr="00D47976" w:rsidRDefault="00BB2E43" w:rsidP="00BB2E43"><w:pPr><w:ind w:left="-85"/><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:b/><w:sz w:val="24"/></w:rPr></w:pPr><w:r w:rsidRPr="00D47976"><w:rPr><w:b/><w:sz w:val="24"/></w:rPr><w:lastRenderedPageBreak/><w:t>Alternate Scenarios</w:t></w:r></w:p>
'@
Set-ParagraphXml "Alternate Scenarios" $xml8
